# Update "想去人数" (want-to-go count) figures on both the "展览" sheet
# and the "全部类型" summary sheet, which mirrors the same rows.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> hashtable of cell -> new value
$updates = @{
    "展览" = @{
        "F2"  = 14081
        "F4"  = 548
        "F8"  = 13932
        "F9"  = 14959
        "F20" = 62
        "F23" = 122
        "F25" = 5816
        "F26" = 947
        "F28" = 5446
        "F32" = 326
    }
    "全部类型" = @{
        "F2"  = 14081
        "F5"  = 548
        "F9"  = 13932
        "F10" = 14959
        "F21" = 62
        "F24" = 122
        "F27" = 5816
        "F28" = 947
        "F30" = 5446
        "F34" = 326
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellUpdates = $updates[$sheetName]
    foreach ($cellRef in $cellUpdates.Keys) {
        $ws.Range($cellRef).Value = $cellUpdates[$cellRef]
    }
}
